$d = $word.ActiveDocument

# 1) Email address: the three runs "g.piacentino" + "@" + "columbia.edu"
#    become a single run "giorgiapiacentino@gmail.com".
#    Replace the first part's text, then delete the now-redundant
#    "@" and "columbia.edu" runs via Find/Replace on the combined text.
$r1 = $d.Content.Find.Execute("g.piacentino@columbia.edu", $true, $true, $false, $false, $false, $true, 1, $false, "giorgiapiacentino@gmail.com", 2)

# 2) Remove ", University of Minnesota (scheduled)" trailing text that
#    followed "... London Business School (scheduled)".
$r2 = $d.Content.Find.Execute(", University of Minnesota (scheduled)", $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
